$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 "blog" slots rotate forward by one series number, freeing up ser:127
# and introducing ser:130. The J8 "signin" card is swapped from the buggy
# double-escaped-quote variant to the correct single-escaped-quote variant.

$c8 = @'
type: blog
width: 2
height: 1
ser: 130
'@

$e8 = @'
type: blog
width: 2
height: 1
ser: 129
'@

$i8 = @'
type: blog
width: 2
height: 1
ser: 128
'@

$j8 = @'
type: signin
width: 2
height: 1
h3.w-half: Sign up to get unlimited access to the entire content of zakatlists
button.primary: Sign In*goto("/signin/home")
button.secondary: Sign Up for Rs 300 / Month*goto("/signup")
'@

$ws.Range("C8").Value = $c8
$ws.Range("E8").Value = $e8
$ws.Range("I8").Value = $i8
$ws.Range("J8").Value = $j8

# Move the active selection from I8 to J8 to match the saved view state.
[void]$ws.Range("J8").Select()
